$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the recording date ("27.01.2025") into A8 and A9: these rows had no
# acceleration value for the given velocity, so the date is recorded
# instead.
$ws.Range("A8").Value = "27.01.2025"
$ws.Range("A9").Value = "27.01.2025"
